# Aprimorando sistema para possibilitar concluir cards
# Adds a new "Baixa" (Tag) column F to the cards sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: width + header + values for the existing data rows (2 and 4,
# matching the source data exactly).
$ws.Columns.Item(6).ColumnWidth = 26.451822916666668

$ws.Range("F1").Value = "Baixa"
$ws.Range("F2").Value = "Teste"
$ws.Range("F4").Value = "Teste"

# New empty, underlined cell at F7 (sibling of the pre-existing G7 cell),
# and move the active selection there.
$ws.Range("F7").Font.Underline = $true
$ws.Range("F7").Select() | Out-Null
